$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New time-tracking entry in row 25: "Testdaten und Entwurf DB erstellen" ---
# Copy the date formatting (short-date number format) from the previous
# populated row (A24) onto A25 before writing the actual date value, so the
# cell picks up the existing date style instead of "General".
[void]$ws.Range("A24").Copy()
[void]$ws.Range("A25").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A25").Value = 45869               # 31.07.2025
$ws.Range("B25").Value = "Testdaten und Entwurf DB erstellen"
$ws.Range("G25").Value = 0.75

# --- New time-tracking entry in row 26: "Arbeitspakete definiert" ---
[void]$ws.Range("A24").Copy()
[void]$ws.Range("A26").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A26").Value = 45718               # 02.03.2025
$ws.Range("B26").Value = "Arbeitspakete definiert"
$ws.Range("G26").Value = 2

$excel.CutCopyMode = 0

# Move the selection to reflect where the user ended up after filling in
# the new rows (the saved file's selection spans B26:F26 and B32:F32, with
# B32 being the active cell - the last selected block is restored here).
[void]$ws.Range("B32:F32").Select()
